$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 11 data: B11 = 5
$ws.Range("B11").Value = 5

# Update selection to F5
$ws.Range("F5").Select()
